$wb = $excel.ActiveWorkbook

# --- HOME sheet: new "Section 2" block (rows 8-13) ---
$home = $wb.Worksheets.Item("HOME")

$home.Cells.Item(8, 1).Value = "SECTION_2_TEXT_1"
$home.Cells.Item(8, 2).Value = "we specialized in"

$home.Cells.Item(9, 1).Value  = "SECTION_2_MENU_1"
$home.Cells.Item(10, 1).Value = "SECTION_2_MENU_2"
$home.Cells.Item(11, 1).Value = "SECTION_2_MENU_3"
$home.Cells.Item(12, 1).Value = "SECTION_2_MENU_4"
$home.Cells.Item(13, 1).Value = "SECTION_2_MENU_5"

$home.Cells.Item(9, 2).Value  = "Brand identity"
$home.Cells.Item(10, 2).Value = "Motion"
$home.Cells.Item(11, 2).Value = "User Interface"
$home.Cells.Item(12, 2).Value = "2D Graphic"
$home.Cells.Item(13, 2).Value = "Digital Project"

# --- LINK sheet: new WORKS_* menu links (rows 14-19) ---
$link = $wb.Worksheets.Item("LINK")

$link.Cells.Item(14, 1).Value = "WORKS_ALL"
$link.Cells.Item(15, 1).Value = "WORKS_BRAND"
$link.Cells.Item(16, 1).Value = "WORKS_MOTION"
$link.Cells.Item(17, 1).Value = "WORKS_INTERFACE"
$link.Cells.Item(18, 1).Value = "WORKS_GRAPHIC"
$link.Cells.Item(19, 1).Value = "WORKS_DIGITAL"

$link.Cells.Item(14, 2).Value = "/works/all"
$link.Cells.Item(15, 2).Value = "/works/brand"
$link.Cells.Item(16, 2).Value = "/works/motion"
$link.Cells.Item(17, 2).Value = "/works/interface"
$link.Cells.Item(18, 2).Value = "/works/graphic"
$link.Cells.Item(19, 2).Value = "/works/digital"

# Widen LINK columns A/B so the new, longer values fit (mirrors Excel's
# automatic best-fit column sizing after the paste).
$link.Columns.Item(1).ColumnWidth = 17.666666666666668
$link.Columns.Item(2).ColumnWidth = 25.833333333333332

# --- Selection / active-sheet bookkeeping ---
$home.Range("B17").Select() | Out-Null
$link.Range("B20").Select() | Out-Null

$link.Activate()
